# Updates cryptos list prices (column D) and 1h volume % (column E)
# Values are forced to Text format ("@") so numeric-looking strings
# (e.g. "603.41") are preserved as text rather than being converted to
# floating point numbers by Excel, matching the original inlineStr data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.410.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.190.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.190.41"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.712.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.463.10"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.192.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.90"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.91"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.65"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.01"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "510.71"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.75"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0424"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.85%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.71%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.91%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.852.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.44"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.34%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.22%  "
